# Commit: test updates because of change to tax calc
# The "raw data" sheet stores cached model outputs as plain numbers (no
# formulas), so the updated tax-calc results are written directly as values,
# matching how the authoring tool regenerated this fixture.

$wb = $excel.ActiveWorkbook
$wsReport = $wb.Worksheets.Item("report")
$wsRaw = $wb.Worksheets.Item("raw data")

# --- report sheet updates ---
$wsReport.Range("B3").Value = 45621
$wsReport.Range("B11").Value = 4539.888998736316
$wsReport.Range("B17").Value = 246.6777894289741

# --- raw data sheet updates (tax calc change) ---
$wsRaw.Range("D7").Value = 81.2351494668635
$wsRaw.Range("C9").Value = 1236.364452457454
$wsRaw.Range("D10").Value = 80.3636894097345
$wsRaw.Range("C11").Value = [double]"-1.4210854715202e-13"
$wsRaw.Range("D11").Value = 333.508252665948
$wsRaw.Range("E11").Value = 771.6510137982533
$wsRaw.Range("F11").Value = 1248.773671883719
$wsRaw.Range("G11").Value = 1767.833752449617
$wsRaw.Range("H11").Value = 2308.532275552996
$wsRaw.Range("I11").Value = 2864.052381650803
$wsRaw.Range("J11").Value = 3426.717477683967
$wsRaw.Range("K11").Value = 3988.182342556486
$wsRaw.Range("L11").Value = 4539.888998736316
$wsRaw.Range("D14").Value = 35.42421237766585
$wsRaw.Range("E14").Value = 29.7739312022838
$wsRaw.Range("F14").Value = 25.05478373173997
$wsRaw.Range("G14").Value = 21.10677593009354
$wsRaw.Range("H14").Value = 17.88920337191616
$wsRaw.Range("I14").Value = 15.25632516352358
$wsRaw.Range("J14").Value = 13.09355009711794
$wsRaw.Range("K14").Value = 11.31067968283281
$wsRaw.Range("L14").Value = 11.31067968283281
$wsRaw.Range("C15").Value = 59.39819546145245
$wsRaw.Range("D15").Value = 77.54353235397711
$wsRaw.Range("E15").Value = 101.1035703087453
$wsRaw.Range("F15").Value = 131.473924734757
$wsRaw.Range("G15").Value = 168.0454794187028
$wsRaw.Range("H15").Value = 210.993468792832
$wsRaw.Range("I15").Value = 260.1588601482575
$wsRaw.Range("J15").Value = 314.9218588035449
$wsRaw.Range("K15").Value = 374.1223234412789
$wsRaw.Range("L15").Value = 374.1223234412789
$wsRaw.Range("C17").Value = [double]"-1.4210854715202e-13"
$wsRaw.Range("D17").Value = 333.5082526659482
$wsRaw.Range("D19").Value = 386.8340450803025
$wsRaw.Range("C20").Value = -416.6069075425462
$wsRaw.Range("D20").Value = 51.92254754254623
$wsRaw.Range("C21").Value = 38.56499405745386
$wsRaw.Range("D21").Value = 81.2351494668635
$wsRaw.Range("D22").Value = 386.8340450803025
$wsRaw.Range("C23").Value = 256.8869075425461
$wsRaw.Range("D23").Value = 281.585705123402
$wsRaw.Range("C24").Value = -159.7200000000001
$wsRaw.Range("D24").Value = 333.5082526659482
$wsRaw.Range("C25").Value = 341.766986878546
$wsRaw.Range("D25").Value = 345.0730197570922
$wsRaw.Range("L27").Value = 48.7774980505551
$wsRaw.Range("B29").Value = 4973.367662396239
$wsRaw.Range("C29").Value = 5630.424428635863
$wsRaw.Range("B30").Value = 6626.33902239624
$wsRaw.Range("C30").Value = 6866.788881093317
$wsRaw.Range("B31").Value = 0.0878640292336085
$wsRaw.Range("C31").Value = 0.09124057377420373
$wsRaw.Range("B32").Value = 6641.359748607475
$wsRaw.Range("C32").Value = 6883.129388831488
$wsRaw.Range("B33").Value = 245.3166451515794
$wsRaw.Range("C33").Value = 269.8483096667374
$wsRaw.Range("D33").Value = 296.8331406334111
$wsRaw.Range("E33").Value = 326.5164546967522
$wsRaw.Range("F33").Value = 359.1681001664275
$wsRaw.Range("G33").Value = 395.0849101830703
$wsRaw.Range("H33").Value = 434.5934012013773
$wsRaw.Range("I33").Value = 478.0527413215151
$wsRaw.Range("J33").Value = 525.8580154536667
$wsRaw.Range("K33").Value = 578.4438169990334
$wsRaw.Range("L33").Value = 587.5107006483817
$wsRaw.Range("B34").Value = 126.4394274323507
$wsRaw.Range("C34").Value = 138.5580649780907
$wsRaw.Range("B35").Value = 246.6777894289741
$wsRaw.Range("C35").Value = 271.209453944132
$wsRaw.Range("D35").Value = 298.4100281984602
$wsRaw.Range("E35").Value = 328.3925925712065
$wsRaw.Range("F35").Value = 361.3976145217733
$wsRaw.Range("G35").Value = 397.7314536525708
$wsRaw.Range("H35").Value = 437.7159550031689
$wsRaw.Range("I35").Value = 481.714173507921
$wsRaw.Range("J35").Value = 530.1242381041834
$wsRaw.Range("K35").Value = 583.3825122469223
$wsRaw.Range("L35").Value = 592.4493958962705
